{"js": "// Update the worksheet date header and the 25 three-digit x one-digit\n// multiplication problems/answers to the new values.\nconst replacements = [\n  [\"2025-09-01 Monday\", \"2025-09-02 Tuesday\"],\n  [\"557\u00d77=3899\", \"703\u00d76=4218\"],\n  [\"183\u00d77=1281\", \"333\u00d72=666\"],\n  [\"505\u00d79=4545\", \"925\u00d74=3700\"],\n  [\"776\u00d79=6984\", \"958\u00d78=7664\"],\n  [\"512\u00d73=1536\", \"231\u00d76=1386\"],\n  [\"750\u00d75=3750\", \"530\u00d72=1060\"],\n  [\"989\u00d75=4945\", \"686\u00d74=2744\"],\n  [\"129\u00d79=1161\", \"353\u00d78=2824\"],\n  [\"391\u00d77=2737\", \"371\u00d74=1484\"],\n  [\"132\u00d79=1188\", \"239\u00d73=717\"],\n  [\"351\u00d73=1053\", \"830\u00d74=3320\"],\n  [\"638\u00d74=2552\", \"183\u00d72=366\"],\n  [\"946\u00d73=2838\", \"747\u00d78=5976\"],\n  [\"444\u00d78=3552\", \"601\u00d77=4207\"],\n  [\"150\u00d79=1350\", \"661\u00d79=5949\"],\n  [\"424\u00d73=1272\", \"837\u00d79=7533\"],\n  [\"316\u00d76=1896\", \"507\u00d79=4563\"],\n  [\"361\u00d78=2888\", \"349\u00d78=2792\"],\n  [\"257\u00d74=1028\", \"773\u00d73=2319\"],\n  [\"924\u00d72=1848\", \"172\u00d76=1032\"],\n  [\"381\u00d79=3429\", \"702\u00d72=1404\"],\n  [\"254\u00d73=762\", \"172\u00d74=688\"],\n  [\"807\u00d79=7263\", \"373\u00d72=746\"],\n  [\"914\u00d72=1828\", \"380\u00d77=2660\"],\n  [\"152\u00d75=760\", \"736\u00d76=4416\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date header and the 25 three-digit x one-digit\n# multiplication problems/answers to the new values.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2025-09-01 Monday\", \"2025-09-02 Tuesday\"),\n  @(\"557\u00d77=3899\", \"703\u00d76=4218\"),\n  @(\"183\u00d77=1281\", \"333\u00d72=666\"),\n  @(\"505\u00d79=4545\", \"925\u00d74=3700\"),\n  @(\"776\u00d79=6984\", \"958\u00d78=7664\"),\n  @(\"512\u00d73=1536\", \"231\u00d76=1386\"),\n  @(\"750\u00d75=3750\", \"530\u00d72=1060\"),\n  @(\"989\u00d75=4945\", \"686\u00d74=2744\"),\n  @(\"129\u00d79=1161\", \"353\u00d78=2824\"),\n  @(\"391\u00d77=2737\", \"371\u00d74=1484\"),\n  @(\"132\u00d79=1188\", \"239\u00d73=717\"),\n  @(\"351\u00d73=1053\", \"830\u00d74=3320\"),\n  @(\"638\u00d74=2552\", \"183\u00d72=366\"),\n  @(\"946\u00d73=2838\", \"747\u00d78=5976\"),\n  @(\"444\u00d78=3552\", \"601\u00d77=4207\"),\n  @(\"150\u00d79=1350\", \"661\u00d79=5949\"),\n  @(\"424\u00d73=1272\", \"837\u00d79=7533\"),\n  @(\"316\u00d76=1896\", \"507\u00d79=4563\"),\n  @(\"361\u00d78=2888\", \"349\u00d78=2792\"),\n  @(\"257\u00d74=1028\", \"773\u00d73=2319\"),\n  @(\"924\u00d72=1848\", \"172\u00d76=1032\"),\n  @(\"381\u00d79=3429\", \"702\u00d72=1404\"),\n  @(\"254\u00d73=762\", \"172\u00d74=688\"),\n  @(\"807\u00d79=7263\", \"373\u00d72=746\"),\n  @(\"914\u00d72=1828\", \"380\u00d77=2660\"),\n  @(\"152\u00d75=760\", \"736\u00d76=4416\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $range = $d.Content\n  $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $found) {\n    throw \"Text not found: $oldText\"\n  }\n}\n"}
